# Weekly update: insert the new week's price data (date 44610) at the top of
# the data block (rows 703-704) for Vega Monumental Concepción - Manzana.
# All existing weekly groups shift down by two rows (rows 703-738 ->
# 705-740); no data is lost, the table simply grows by one week.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at the top of the data block; everything from
# row 703 down shifts down by 2 rows (703-738 -> 705-740).
$ws.Rows("703:704").Insert()

# Fill in the two new rows with this week's data.
$newRows = @(
    @(11, "Vega Monumental Concepción", "Bíobío", 44610, 8, "Fruta", 100104, "Frutos de pepita", 100104002, "Manzana", "Granny Smith", "Segunda", 130, 10000, 11000, 10615, "`$/caja 16 kilos empedrada", "Región de O'Higgins", 663, 16),
    @(11, "Vega Monumental Concepción", "Bíobío", 44610, 8, "Fruta", 100104, "Frutos de pepita", 100104002, "Manzana", "Royal Gala",   "Primera", 180, 11000, 12000, 11444, "`$/caja 16 kilos empedrada", "Región de O'Higgins", 715, 16)
)

$startRow = 703
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $newRows[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le $row.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}
